$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: USUARIO "ivan" -> "Manrique", CONTRASEÑA 1234 -> "Nahia82" (text)
$ws.Range("A2").Value = "Manrique"
$ws.Range("B2").Value = "Nahia82"

# Row 3: CONTRASEÑA 1234 -> "Aa123456" (text)
$ws.Range("B3").Value = "Aa123456"

# Row 4: CONTRASEÑA 1234 -> 25 (number)
$ws.Range("B4").Value = 25

# Update the active selection to D6 (matches the saved view state in the diff)
$ws.Range("D6").Select()
